# Scheduled runner update: refresh market-price-derived Leve profit figures
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3240.8
$ws.Range("I38").Value = 10.75
$ws.Range("J38").Value = 5394.1665
$ws.Range("K38").Value = 32.25
$ws.Range("L38").Value = 16182.4995
$ws.Range("M38").Value = 339.75
$ws.Range("N38").Value = -16926.4995
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 10000
$ws.Range("K113").Value = 10000
$ws.Range("M113").Value = -6746
$ws.Range("H132").Value = 4288.6665
$ws.Range("J132").Value = 9000
$ws.Range("L132").Value = 27000
$ws.Range("N132").Value = -32060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1498.75
$ws.Range("I32").Value = 1498.75
$ws.Range("K32").Value = 1498.75
$ws.Range("M32").Value = -1211.75
$ws.Range("H61").Value = 10232.833
$ws.Range("I61").Value = 5599.25
$ws.Range("K61").Value = 5599.25
$ws.Range("M61").Value = -5387.25
$ws.Range("H74").Value = 9900
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 21000
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 21000
$ws.Range("M74").Value = -1626
$ws.Range("N74").Value = -22748
$ws.Range("H77").Value = 9900
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 21000
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 105000
$ws.Range("M77").Value = -8132
$ws.Range("H102").Value = 1183
$ws.Range("I102").Value = 1183
$ws.Range("K102").Value = 1183
$ws.Range("M102").Value = 439
$ws.Range("H132").Value = 17809.182
$ws.Range("I132").Value = 15989
$ws.Range("J132").Value = 19993.4
$ws.Range("K132").Value = 47967
$ws.Range("L132").Value = 59980.2
$ws.Range("M132").Value = -45437
$ws.Range("N132").Value = -65040.2
$ws.Range("H136").Value = 10232.833
$ws.Range("I136").Value = 5599.25
$ws.Range("K136").Value = 16797.75
$ws.Range("M136").Value = -14247.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1956.1428
$ws.Range("I105").Value = 1658.8
$ws.Range("K105").Value = 1658.8
$ws.Range("M105").Value = 88.20000000000005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H6").Value = 14000000
$ws.Range("I6").Value = 14000000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 14000000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -13999887
$ws.Range("N6").ClearContents()
$ws.Range("H31").Value = 18333.334
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 18333.334
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 18333.334
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -18923.334
$ws.Range("H34").Value = 18333.334
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 18333.334
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 18333.334
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -18737.334
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H99").Value = 5666.6665
$ws.Range("I99").Value = 5500
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 5500
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -4002
$ws.Range("N99").Value = -8996
$ws.Range("H126").Value = 5666.6665
$ws.Range("I126").Value = 5500
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 16500
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -14030
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 144.29033
$ws.Range("I2").Value = 170.34616
$ws.Range("J2").Value = 8.800000000000001
$ws.Range("K2").Value = 1022.07696
$ws.Range("L2").Value = 52.8
$ws.Range("M2").Value = -909.07696
$ws.Range("N2").Value = -278.8
$ws.Range("H38").Value = 67.333336
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H59").Value = 2552.5
$ws.Range("I59").Value = 105
$ws.Range("K59").Value = 315
$ws.Range("M59").Value = 225
$ws.Range("H140").Value = 1187
$ws.Range("I140").Value = 1187
$ws.Range("K140").Value = 3561
$ws.Range("M140").Value = 1619

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3699.3333
$ws.Range("I80").Value = 4100
$ws.Range("K80").Value = 4100
$ws.Range("M80").Value = -3102
$ws.Range("H83").Value = 3699.3333
$ws.Range("I83").Value = 4100
$ws.Range("K83").Value = 20500
$ws.Range("M83").Value = -15508
$ws.Range("H138").Value = 89997.5
$ws.Range("J138").Value = 89997.5
$ws.Range("L138").Value = 89997.5
$ws.Range("N138").Value = -100277.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 950
$ws.Range("I22").Value = 950
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 950
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -655
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 950
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -843
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 4287.375
$ws.Range("I46").Value = 4159.8
$ws.Range("K46").Value = 4159.8
$ws.Range("M46").Value = -3971.8
$ws.Range("H55").Value = 852.5714
$ws.Range("J55").Value = 496.66666
$ws.Range("L55").Value = 496.66666
$ws.Range("N55").Value = -842.66666
$ws.Range("H68").Value = 4118.5
$ws.Range("I68").Value = 4495
$ws.Range("J68").Value = 3993
$ws.Range("K68").Value = 4495
$ws.Range("L68").Value = 3993
$ws.Range("M68").Value = -3746
$ws.Range("N68").Value = -5491
$ws.Range("H71").Value = 4118.5
$ws.Range("I71").Value = 4495
$ws.Range("J71").Value = 3993
$ws.Range("K71").Value = 22475
$ws.Range("L71").Value = 19965
$ws.Range("M71").Value = -18731
$ws.Range("N71").Value = -27453
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2500
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2500
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H104").Value = 36600
$ws.Range("J104").Value = 36600
$ws.Range("L104").Value = 36600
$ws.Range("N104").Value = -43588
$ws.Range("H122").Value = 2499
$ws.Range("I122").Value = 2497
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 7491
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -5041
$ws.Range("N122").Value = -12400
